# Changes to accommodate multiple market places
# - Add new integration columns (Trendyol discount fields, Amazon(v3) fields) to the CEOrder sheet
# - Switch the active/selected tab from "Retailer" to "CEOrder"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CEOrder")

# --- New header cells (row 1) ---
# Order chosen to match the original authoring order of new shared strings.
$ws.Range("N1").Value = "Amazon(v3).intAE_forwardCarriyoBooking"
$ws.Range("L1").Value = "Trendyol.intAE_lineLevelDiscountField"
$ws.Range("M1").Value = "Trendyol.intAE_orderLevelDiscountField"
$ws.Range("L2").Value = "UNIT_DISCOUNT_AMOUNT"
$ws.Range("M2").Value = "TOTAL_DISCOUNT_AMOUNT"
$ws.Range("O1").Value = "Amazon(v3).intAE_lineLevelDiscountField"
$ws.Range("P1").Value = "Amazon(v3).intAE_orderLevelDiscountField"

# --- New value cells (row 2) ---
$ws.Range("N2").Value = $true
$ws.Range("O2").Value = "UNIT_DISCOUNT_AMOUNT"
$ws.Range("P2").Value = "TOTAL_DISCOUNT_AMOUNT"

# --- Column widths for the newly populated columns ---
$ws.Columns.Item(6).ColumnWidth = 7.1
$ws.Columns.Item(7).ColumnWidth = 9.3
$ws.Columns.Item(8).ColumnWidth = 33.65
$ws.Columns.Item(9).ColumnWidth = 32.1
$ws.Columns.Item(11).ColumnWidth = 38.45
$ws.Columns.Item(12).ColumnWidth = 38.45
$ws.Columns.Item(13).ColumnWidth = 38.45

# --- Make CEOrder the active/selected sheet (this also clears tabSelected on
#     whichever sheet was previously active, i.e. "Retailer") and move the
#     selection to P1 to match the new last populated header cell ---
$ws.Activate()
$ws.Range("P1").Select()
